$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Row 54
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("L54").ClearContents()
$ws.Range("N54").Value = 0

# Row 55
$ws.Range("H55").Value = 440.9
$ws.Range("I55").Value = 68.166664
$ws.Range("J55").Value = 1000
$ws.Range("K55").Value = 68.166664
$ws.Range("L55").Value = 1000
$ws.Range("M55").Value = 145.833336
$ws.Range("N55").Value = -1428

# Row 94
$ws.Range("H94").Value = 4994.75
$ws.Range("I94").Value = 4994.75
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 4994.75
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -4543.75

# Row 107
$ws.Range("H107").Value = 283.8
$ws.Range("I107").Value = 232.375
$ws.Range("J107").Value = 489.5
$ws.Range("K107").Value = 232.375
$ws.Range("L107").Value = 489.5
$ws.Range("M107").Value = 1687.625
$ws.Range("N107").Value = -4329.5

# Row 135
$ws.Range("H135").Value = 2660.8572
$ws.Range("I135").Value = 2131.75
$ws.Range("J135").Value = 3366.3333
$ws.Range("K135").Value = 19185.75
$ws.Range("L135").Value = 30296.9997
$ws.Range("M135").Value = -16650.75
$ws.Range("N135").Value = -35366.9997

# Row 138
$ws.Range("H138").Value = 4221.7144
$ws.Range("I138").Value = 7998.5
$ws.Range("J138").Value = 3592.25
$ws.Range("K138").Value = 23995.5
$ws.Range("L138").Value = 10776.75
$ws.Range("M138").Value = -18855.5
$ws.Range("N138").Value = -21056.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Range("H32").Value = 4180.387
$ws.Range("I32").Value = 4180.387
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 4180.387
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -3893.387

# Row 50
$ws.Range("H50").Value = 44000
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 44000
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 44000
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -45428

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

# Row 134
$ws.Range("H134").Value = 1956.2941
$ws.Range("I134").Value = 1580
$ws.Range("J134").Value = 4778.5
$ws.Range("K134").Value = 4740
$ws.Range("L134").Value = 14335.5
$ws.Range("M134").Value = -2205
$ws.Range("N134").Value = -19405.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 4231.4287
$ws.Range("I31").Value = 2929.818
$ws.Range("J31").Value = 9004
$ws.Range("K31").Value = 2929.818
$ws.Range("L31").Value = 9004
$ws.Range("M31").Value = -2634.818
$ws.Range("N31").Value = -9594

# Row 34
$ws.Range("H34").Value = 4231.4287
$ws.Range("I34").Value = 2929.818
$ws.Range("J34").Value = 9004
$ws.Range("K34").Value = 2929.818
$ws.Range("L34").Value = 9004
$ws.Range("M34").Value = -2727.818
$ws.Range("N34").Value = -9408

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

# Row 63
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").ClearContents()
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = 0

# Row 66
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").ClearContents()
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = 0

# Row 131
$ws.Range("H131").Value = 1477
$ws.Range("I131").Value = 1000
$ws.Range("J131").Value = 1715.5
$ws.Range("K131").Value = 3000
$ws.Range("L131").Value = 5146.5
$ws.Range("M131").Value = 2040
$ws.Range("N131").Value = -15226.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

# Row 69
$ws.Range("H69").Value = 10000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 10000
$ws.Range("K69").Value = 0
$ws.Range("L69").ClearContents()
$ws.Range("M69").Value = 10000
$ws.Range("N69").Value = -11498

# Row 72
$ws.Range("H72").Value = 10000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 10000
$ws.Range("K72").Value = 0
$ws.Range("L72").ClearContents()
$ws.Range("M72").Value = 30000
$ws.Range("N72").Value = -37488

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

# Row 4
$ws.Range("H4").Value = 4000
$ws.Range("I4").Value = 4000
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 4000
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -3887

# Row 17
$ws.Range("H17").Value = 18500
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 18500
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 18500
$ws.Range("N17").Value = -18840

# Row 26
$ws.Range("H26").Value = 3336.3333
$ws.Range("I26").Value = 9
$ws.Range("J26").Value = 5000
$ws.Range("K26").Value = 9
$ws.Range("L26").Value = 5000
$ws.Range("M26").Value = 286
$ws.Range("N26").Value = -5590

# Row 28
$ws.Range("H28").Value = 4000
$ws.Range("I28").Value = 4000
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 4000
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -3768

# Row 29
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()

# Row 31
$ws.Range("H31").Value = 23008004
$ws.Range("I31").Value = 4015
$ws.Range("J31").Value = 34510000
$ws.Range("K31").Value = 4015
$ws.Range("L31").Value = 34510000
$ws.Range("M31").Value = -3767
$ws.Range("N31").Value = -34510496

# Row 34
$ws.Range("H34").Value = 16332.667
$ws.Range("I34").Value = 14500
$ws.Range("J34").Value = 19998
$ws.Range("K34").Value = 14500
$ws.Range("L34").Value = 19998
$ws.Range("M34").Value = -14328
$ws.Range("N34").Value = -20342

# Row 37
$ws.Range("H37").Value = 4000
$ws.Range("I37").Value = 4000
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 4000
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -3893

# Row 41
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()

# Row 53
$ws.Range("H53").Value = 33681.332
$ws.Range("I53").Value = 33681.332
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 33681.332
$ws.Range("L53").ClearContents()
$ws.Range("N53").Value = 0
$ws.Range("M53").Value = -33163.332

# Row 58
$ws.Range("H58").Value = 13103
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 13103
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 13103
$ws.Range("N58").Value = -13623

# Row 93
$ws.Range("H93").Value = 1070.6666
$ws.Range("I93").Value = 885.1
$ws.Range("J93").Value = 1302.625
$ws.Range("K93").Value = 885.1
$ws.Range("L93").Value = 1302.625
$ws.Range("M93").Value = 362.9
$ws.Range("N93").Value = -3798.625

# Row 122
$ws.Range("H122").Value = 5897.3335
$ws.Range("I122").Value = 5278.8
$ws.Range("J122").Value = 8990
$ws.Range("K122").Value = 15836.4
$ws.Range("L122").Value = 26970
$ws.Range("M122").Value = -13386.4
$ws.Range("N122").Value = -31870

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

# Row 19
$ws.Range("H19").Value = 5000
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 5000
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 5000
$ws.Range("N19").Value = -5348

# Row 32
$ws.Range("H32").Value = 2000
$ws.Range("I32").Value = 2000
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -1683

# Row 40
$ws.Range("H40").Value = 20000
$ws.Range("I40").Value = 20000
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 20000
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -19851

# Row 44
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").ClearContents()
$ws.Range("N44").Value = 0

# Row 62
$ws.Range("H62").Value = 9999.666999999999
$ws.Range("I62").Value = 9999.666999999999
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 9999.666999999999
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -9375.666999999999

# Row 65
$ws.Range("H65").Value = 9999.666999999999
$ws.Range("I65").Value = 9999.666999999999
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 49998.335
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -46878.335

# Row 81
$ws.Range("H81").Value = 725
$ws.Range("I81").Value = 700
$ws.Range("J81").Value = 750
$ws.Range("K81").Value = 1400
$ws.Range("L81").Value = 1500
$ws.Range("M81").Value = -339
$ws.Range("N81").Value = -3622

# Row 84
$ws.Range("H84").Value = 725
$ws.Range("I84").Value = 700
$ws.Range("J84").Value = 750
$ws.Range("K84").Value = 7000
$ws.Range("L84").Value = 7500
$ws.Range("M84").Value = -1696
$ws.Range("N84").Value = -18108
